$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1750663129973475
$ws.Range("C2").Value = 0.5570291777188329
$ws.Range("J2").Value = 0.01856763925729443
$ws.Range("P2").Value = 0.1246684350132626
$ws.Range("S2").Value = 0.1246684350132626
$ws.Range("B3").Value = 0.008733624454148471
$ws.Range("C3").Value = 0.01746724890829694
$ws.Range("J3").Value = 0.03930131004366812
$ws.Range("P3").Value = 0.777292576419214
$ws.Range("S3").Value = 0.1572052401746725
$ws.Range("J4").Value = 0.02
$ws.Range("P4").Value = 0.78
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.08494208494208494
$ws.Range("D6").Value = 0.01158301158301158
$ws.Range("F6").Value = 0.05019305019305019
$ws.Range("J6").Value = 0.2432432432432433
$ws.Range("O6").Value = 0.02702702702702703
$ws.Range("Q6").Value = 0.1776061776061776
$ws.Range("R6").Value = 0.08494208494208494
$ws.Range("S6").Value = 0.3204633204633205
$ws.Range("B7").Value = 0.1068376068376068
$ws.Range("D7").Value = 0.008547008547008548
$ws.Range("F7").Value = 0.05982905982905983
$ws.Range("J7").Value = 0.1581196581196581
$ws.Range("O7").Value = 0.0170940170940171
$ws.Range("Q7").Value = 0.1623931623931624
$ws.Range("R7").Value = 0.1068376068376068
$ws.Range("S7").Value = 0.3803418803418803
$ws.Range("B8").Value = 0.1058315334773218
$ws.Range("D8").Value = 0.02159827213822894
$ws.Range("F8").Value = 0.06263498920086392
$ws.Range("J8").Value = 0.09287257019438445
$ws.Range("O8").Value = 0.01943844492440605
$ws.Range("Q8").Value = 0.1900647948164147
$ws.Range("R8").Value = 0.1360691144708423
$ws.Range("S8").Value = 0.3714902807775378
$ws.Range("B9").Value = 0.08928571428571429
$ws.Range("D9").Value = 0.01785714285714286
$ws.Range("F9").Value = 0.07142857142857142
$ws.Range("J9").Value = 0.08482142857142858
$ws.Range("O9").Value = 0.008928571428571428
$ws.Range("Q9").Value = 0.1875
$ws.Range("R9").Value = 0.125
$ws.Range("S9").Value = 0.4151785714285715
$ws.Range("B10").Value = 0.1109725685785536
$ws.Range("D10").Value = 0.01932668329177057
$ws.Range("E10").Value = 0.001246882793017456
$ws.Range("F10").Value = 0.06795511221945137
$ws.Range("J10").Value = 0.1309226932668329
$ws.Range("O10").Value = 0.01433915211970075
$ws.Range("Q10").Value = 0.2219451371571072
$ws.Range("R10").Value = 0.1034912718204489
$ws.Range("S10").Value = 0.3298004987531172
$ws.Range("G11").Value = 0.1375661375661376
$ws.Range("J11").Value = 0.08201058201058201
$ws.Range("K11").Value = 0.208994708994709
$ws.Range("L11").Value = 0.5502645502645502
$ws.Range("S11").Value = 0.02116402116402116
$ws.Range("G12").Value = 0.7358490566037735
$ws.Range("J12").Value = 0.2075471698113208
$ws.Range("K12").Value = 0.01415094339622642
$ws.Range("L12").Value = 0.01886792452830189
$ws.Range("S12").Value = 0.02358490566037736
$ws.Range("G13").Value = 0.5471698113207547
$ws.Range("J13").Value = 0.3962264150943396
$ws.Range("S13").Value = 0.05660377358490566
$ws.Range("F15").Value = 0.0158102766798419
$ws.Range("H15").Value = 0.1422924901185771
$ws.Range("I15").Value = 0.07905138339920949
$ws.Range("J15").Value = 0.3952569169960474
$ws.Range("K15").Value = 0.06719367588932806
$ws.Range("M15").Value = 0.003952569169960474
$ws.Range("N15").Value = 0.003952569169960474
$ws.Range("O15").Value = 0.03162055335968379
$ws.Range("S15").Value = 0.2608695652173913
$ws.Range("F16").Value = 0.012
$ws.Range("H16").Value = 0.136
$ws.Range("I16").Value = 0.052
$ws.Range("J16").Value = 0.42
$ws.Range("K16").Value = 0.136
$ws.Range("M16").Value = 0.02
$ws.Range("N16").Value = 0.008
$ws.Range("O16").Value = 0.048
$ws.Range("S16").Value = 0.168
$ws.Range("F17").Value = 0.006993006993006993
$ws.Range("H17").Value = 0.1398601398601399
$ws.Range("I17").Value = 0.1013986013986014
$ws.Range("J17").Value = 0.4527972027972028
$ws.Range("K17").Value = 0.08566433566433566
$ws.Range("M17").Value = 0.01748251748251748
$ws.Range("N17").Value = 0.001748251748251748
$ws.Range("O17").Value = 0.07342657342657342
$ws.Range("S17").Value = 0.1206293706293706
$ws.Range("F18").Value = 0.0126984126984127
$ws.Range("H18").Value = 0.1396825396825397
$ws.Range("I18").Value = 0.06349206349206349
$ws.Range("J18").Value = 0.4507936507936508
$ws.Range("K18").Value = 0.09206349206349207
$ws.Range("M18").Value = 0.0253968253968254
$ws.Range("N18").Value = 0.006349206349206349
$ws.Range("O18").Value = 0.1015873015873016
$ws.Range("S18").Value = 0.1079365079365079
$ws.Range("F19").Value = 0.01819454163750875
$ws.Range("H19").Value = 0.1889433170048985
$ws.Range("I19").Value = 0.07907627711686493
$ws.Range("J19").Value = 0.384184744576627
$ws.Range("K19").Value = 0.1210636808957313
$ws.Range("M19").Value = 0.02169349195241427
$ws.Range("N19").Value = 0.002099370188943317
$ws.Range("O19").Value = 0.06368089573128062
$ws.Range("S19").Value = 0.1210636808957313
